$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 22:05"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1848672
$ws.Range("C4").Value = 11502
$ws.Range("D4").Value = 607323
$ws.Range("E4").Value = 1134807
$ws.Range("G4").Value = 347
$ws.Range("H4").Value = 106542

# Row 10: India
$ws.Range("B10").Value = 198370
$ws.Range("C10").Value = 7761
$ws.Range("D10").Value = 95754
$ws.Range("E10").Value = 97008

# Row 17: Canada
$ws.Range("D17").Value = 49476
$ws.Range("E17").Value = 34846

# Row 38: Irlanda
$ws.Range("E38").Value = 1323
$ws.Range("H38").Value = 1650

# Rows 77/78: Senegal/Guinea swap
$ws.Range("A77").Value = "Guinea"
$ws.Range("B77").Value = 3844
$ws.Range("C77").Value = 138
$ws.Range("D77").Value = 2135
$ws.Range("E77").Value = 1686
$ws.Range("H77").Value = 23

$ws.Range("A78").Value = "Senegal"
$ws.Range("B78").Value = 3739
$ws.Range("C78").Value = 94
$ws.Range("D78").Value = 1858
$ws.Range("E78").Value = 1839
$ws.Range("H78").Value = 42

# Rows 201/202: Santa Lucia/Belice swap
$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("D202").Value = 18
$ws.Range("H202").Value = 0

# Rows 213/214: Papua Nueva Guinea/Islas Virgenes Britanicas swap
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
